$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestSteps")
$ws2 = $wb.Worksheets.Item("TestData")

# --- TestSteps: insert two new rows right after the header row ---
$ws1.Rows("2:3").Insert()

# Copy formatting from the row that used to be row 2 (now shifted to row 4,
# still carrying the original data-row style) into the two freshly inserted
# rows so they match the rest of the table's look (fill/border/font).
$fmtSource = $ws1.Range("A4:F4")
$fmtSource.Copy()
$ws1.Range("A2:F3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New step 1: accessibility check added before the existing flow
$ws1.Range("A2").Value = "smallWaitForElementPresent"
$ws1.Range("B2").Value = "txt_programCode_pcp"
$ws1.Range("C2").Value = "getData=WaitForPageLoad"

# New step 2: the actual axe/accessibility check
$ws1.Range("A3").Value = "checkAccessibility"
$ws1.Range("B3").Value = "AddNewProfile_FundeProfile"

# --- TestSteps: append two blank rows at the bottom (24, 25) ---
$fmtBottom = $ws1.Range("A23:F23")
$fmtBottom.Copy()
$ws1.Range("A24:F25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("A24:C25").ClearContents()

# --- Sheet view / selection / active sheet bookkeeping ---
$ws2.Activate()
$ws2.Range("M4").Select()

$ws1.Activate()
$ws1.Range("D8").Select()
